$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that currently sits after "precedence"
# ------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
    $goBack.Delete()
} catch {
}

# ------------------------------------------------------------------
# 2. Locate the "commands that return a value" table (the one with the
#    int/abs/min/max rows) and update the int example "7.32" -> "7.8"
# ------------------------------------------------------------------
$table = $null
foreach ($t in $d.Tables) {
    for ($i = 1; $i -le $t.Rows.Count; $i++) {
        $c1 = $t.Rows.Item($i).Cells.Item(1).Range.Text
        if ($c1 -like "sqrt*") {
            $table = $t
        }
    }
}

$intRowIndex = 0
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $nameText = $table.Rows.Item($i).Cells.Item(1).Range.Text
    if ($nameText -like "int*") {
        $intRowIndex = $i
    }
}

$intRow = $table.Rows.Item($intRowIndex)
$exampleCell = $intRow.Cells.Item(5)
$find = $exampleCell.Range.Find
$replaceResult = $find.Execute("7.32", $true, $false, $false, $false, $false, $true, 1, $false, "7.8", 2)

# ------------------------------------------------------------------
# 3. Insert a new row right after the int row, for the trunc command
# ------------------------------------------------------------------
$absRowIndex = $intRowIndex + 1
$absRow = $table.Rows.Item($absRowIndex)
$newRow = $table.Rows.Add($absRow)

$newRow.Cells.Item(1).Range.Text = "trunc"
$newRow.Cells.Item(2).Range.Text = "1"
$newRow.Cells.Item(4).Range.Text = "truncates integer part"
$newRow.Cells.Item(5).Range.Text = "trunc 7.8"

# ------------------------------------------------------------------
# 4. Re-create the _GoBack bookmark at the end of the abs row (which
#    now follows the newly inserted trunc row)
# ------------------------------------------------------------------
$absRowIndex2 = $intRowIndex + 2
$absRow2 = $table.Rows.Item($absRowIndex2)
$absExampleCell = $absRow2.Cells.Item(5)
$lastCharRange = $d.Range($absExampleCell.Range.End - 2, $absExampleCell.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $lastCharRange)

# ------------------------------------------------------------------
# 5. Remove the stray paragraph that only contains a page break right
#    after "...for a different output at each frame"
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$found = $find2.Execute("for a different output at each frame", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraRange = $find2.Parent
    $expandResult1 = $paraRange.Expand(4)
    $nextParaRange = $d.Range($paraRange.End, $paraRange.End)
    $expandResult2 = $nextParaRange.Expand(4)
    $nextParaRange.Delete()
}
